# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Refresh the MSME Bulgaria Summary indicator figures with more precise
# (two decimal place) values pulled from the source data.
#
# Values are assigned with a leading apostrophe so Excel keeps them as text
# (matching how the workbook already stores these numeric-looking figures),
# then the style is reset to "Normal" so the quote-prefix indicator doesn't
# leave a stray format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# Enterprises density (per 1000 people) - Statistical Institution source
Set-TextValue "B11" "46.95"
Set-TextValue "C11" "3.88"
Set-TextValue "D11" "50.83"

# Employment (% of total) - Statistical Institution source
Set-TextValue "B12" "31.03"
Set-TextValue "D12" "75.33"

# Enterprises density (per 1000 people) - SME Associations source
Set-TextValue "B33" "40.34"
Set-TextValue "C33" "3.69"
Set-TextValue "D33" "44.03"

# Employment (% of total) - SME Associations source
Set-TextValue "B34" "29.77"
Set-TextValue "C34" "45.44"
Set-TextValue "D34" "75.21"

# Enterprises (% of total) - SME Associations source
Set-TextValue "B36" "91.43"
Set-TextValue "C36" "8.36"
Set-TextValue "D36" "99.79"

# Value added to the economy (% of total) - SME Performance Review EU source
Set-TextValue "B40" "16.49"
Set-TextValue "C40" "43.74"
Set-TextValue "D40" "60.23"
